# Update cryptocurrency price/volume data (and the swapped Bittensor /
# InternetComputer(DFINITY) rows 27-28) to match the latest scrape.
#
# Note: columns B/C/D/E are stored as text in this workbook (e.g. prices use
# '.' as a thousands separator: '63.614.45'). Plain numeric-looking values
# (like '571.44') would otherwise be auto-converted to real numbers by Excel,
# which both changes the cell type and introduces floating point drift
# (e.g. 578.90 -> 578.8999999999999...). Prefixing such values with a leading
# apostrophe forces Excel to keep/store them as literal text, matching the
# original data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.614.45'
$ws.Range('E2').Value = '  -3.32%  '
$ws.Range('D3').Value = '2.607.36'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''571.44'
$ws.Range('E5').Value = '  -4.62%  '
$ws.Range('D6').Value = '''154.98'
$ws.Range('E6').Value = '  -3.09%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -3.35%  '
$ws.Range('D9').Value = '2.603.59'
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('E10').Value = '  -8.21%  '
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '''0.378'
$ws.Range('E12').Value = '  -5.27%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '''27.87'
$ws.Range('E14').Value = '  -4.74%  '
$ws.Range('D15').Value = '3.074.48'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('E16').Value = '  -8.09%  '
$ws.Range('D17').Value = '63.501.32'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '2.615.75'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').Value = '''11.92'
$ws.Range('E19').Value = '  -5.15%  '
$ws.Range('D20').Value = '''7.48'
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('E21').Value = '  -6.52%  '
$ws.Range('D22').Value = '''340.58'
$ws.Range('E22').Value = '  -3.92%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '''67.13'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').Value = '''0.0000106'
$ws.Range('E26').Value = '  -6.22%  '
$ws.Range('B27').Value = 'Bittensor'
$ws.Range('C27').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D27').Value = '''578.90'
$ws.Range('E27').Value = '  +2.55%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''9.05'
$ws.Range('E28').Value = '  -7.45%  '
$ws.Range('E29').Value = '  -4.64%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('E31').Value = '  -2.22%  '
$ws.Range('D32').Value = '''7.80'
$ws.Range('E32').Value = '  -4.26%  '
$ws.Range('E33').Value = '  -4.95%  '
$ws.Range('E34').Value = '  -6.16%  '
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('D36').Value = '''5.31'
$ws.Range('E36').Value = '  -3.57%  '
$ws.Range('E37').Value = '  -5.41%  '
$ws.Range('D38').Value = '''0.999'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  -5.13%  '
$ws.Range('D40').Value = '''153.96'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('E41').Value = '  -6.36%  '
$ws.Range('D43').Value = '''41.55'
$ws.Range('E43').Value = '  -3.18%  '
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('E45').Value = '  -3.05%  '
$ws.Range('D46').Value = '''23.48'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('E47').Value = '  -6.09%  '
$ws.Range('D48').Value = '''0.0581'
$ws.Range('E48').Value = '  -5.96%  '
$ws.Range('D49').Value = '''0.627'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('E50').Value = '  -2.05%  '
$ws.Range('D51').Value = '''0.0245'
$ws.Range('E51').Value = '  -5.17%  '
